$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for the affected rows based on repulled data
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -5
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = -1
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = 0
